$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: I1 = "I0", J1 = "IF" (same style as the other headers)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-11: column I is always 1, column J mirrors column H
for ($r = 2; $r -le 11; $r++) {
    $hValue = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hValue
}
